# "Specular lighting" is done, so remove it from the ToDo list.
# Deleting the entire row shifts every row below it up by one, which is
# exactly what the target workbook shows (rows 3-11 become rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

# Excel leaves the selection on the row that slid up into the deleted
# row's place, now selected as a full row.
$ws.Rows.Item(2).Select()
